$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 is the multiplier vector itself: a disease name in column A,
# followed by the multiplier values (2, 4, 6, ... step 2) across B:BK.
$diseases = @("Anemia", "Asthma", "Cataracts", "Cold", "Cold Sore", "Crohn's Disorder", "Diabetes Type I", "Flu", "HIV", "Meningitis")

$ws.Cells.Item(2, 1).Value = $diseases[0]
for ($col = 2; $col -le 63; $col++) {
    $ws.Cells.Item(2, $col).Value = ($col - 1) * 2
}

# Rows 3-11: just testing some additional disease labels in column A.
for ($i = 1; $i -lt $diseases.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $diseases[$i]
}

# A1 used to hold "N/A" for the (unused) first disease row label;
# rename it to "NaN" to match the multiplier vector header.
$ws.Range("A1").Value = "NaN"

# Move the active selection to reflect where the user was last working.
$ws.Range("H4").Select()
